$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.216.50"
$ws.Range("E2").Value = "  -2.37%  "
$ws.Range("D3").Value = "2.583.01"
$ws.Range("E3").Value = "  -2.55%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'560.27"
$ws.Range("E5").Value = "  -1.89%  "
$ws.Range("D6").Value = "'142.94"
$ws.Range("E6").Value = "  -2.92%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  -1.70%  "
$ws.Range("D9").Value = "2.590.56"
$ws.Range("E9").Value = "  -3.25%  "
$ws.Range("E10").Value = "  -3.23%  "
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("E12").Value = "  +10.63%  "
$ws.Range("E13").Value = "  +3.15%  "
$ws.Range("D14").Value = "3.039.63"
$ws.Range("E14").Value = "  -2.81%  "
$ws.Range("D15").Value = "59.179.46"
$ws.Range("E15").Value = "  -2.37%  "
$ws.Range("D16").Value = "'23.00"
$ws.Range("E16").Value = "  +5.28%  "
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("D18").Value = "2.578.15"
$ws.Range("E18").Value = "  -3.68%  "
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D20").Value = "'336.57"
$ws.Range("E20").Value = "  -2.40%  "
$ws.Range("D21").Value = "'10.35"
$ws.Range("E21").Value = "  -1.25%  "
$ws.Range("D22").Value = "'6.42"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").Value = "'64.01"
$ws.Range("E24").Value = "  -4.07%  "
$ws.Range("E25").Value = "  +5.10%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.81%  "
$ws.Range("E27").Value = "  -2.97%  "
$ws.Range("D28").Value = "'7.34"
$ws.Range("E28").Value = "  -0.86%  "
$ws.Range("D29").Value = "0.0₃0774"
$ws.Range("E29").Value = "  -2.18%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").Value = "'6.14"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("E32").Value = "  -3.20%  "
$ws.Range("D33").Value = "'158.86"
$ws.Range("E33").Value = "  +2.28%  "
$ws.Range("D34").Value = "'19.00"
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("E35").Value = "  -1.94%  "
$ws.Range("D36").Value = "'1.17"
$ws.Range("E36").Value = "  -1.77%  "
$ws.Range("D37").Value = "'0.880"
$ws.Range("E37").Value = "  -3.83%  "
$ws.Range("D38").Value = "'0.868"
$ws.Range("E38").Value = "  -5.33%  "
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("E40").Value = "  -2.74%  "
$ws.Range("D41").Value = "'3.67"
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").Value = "'292.52"
$ws.Range("E42").Value = "  -5.40%  "
$ws.Range("D43").Value = "'132.83"
$ws.Range("E43").Value = "  +4.52%  "
$ws.Range("E45").Value = "  -0.88%  "
$ws.Range("D46").Value = "'0.596"
$ws.Range("E46").Value = "  -2.07%  "
$ws.Range("D47").Value = "'10.63"
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("E48").Value = "  -2.92%  "
$ws.Range("D50").Value = "1.953.64"
$ws.Range("E50").Value = "  -0.76%  "
$ws.Range("D51").Value = "'18.61"
$ws.Range("E51").Value = "  -1.89%  "
